$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "correct" column values from "j" to "left"
$ws.Range("G2").Value = "left"
$ws.Range("G3").Value = "left"

# Reset these cells back to the default "Normal" style, removing the
# center-alignment formatting they previously had
$ws.Range("G2:G3").Style = "Normal"

# Move the active selection to F8, matching the saved view state
$ws.Range("F8").Select()
